$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 16529.166
$ws.Range("I6").Value = 19832
$ws.Range("K6").Value = 59496
$ws.Range("M6").Value = -59384
$ws.Range("H12").Value = 149.66667
$ws.Range("I12").Value = 149.5
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 149.5
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = 20.5
$ws.Range("N12").Value = -490
$ws.Range("H17").Value = 2568.3635
$ws.Range("I17").Value = 1900
$ws.Range("J17").Value = 3125.3333
$ws.Range("K17").Value = 5700
$ws.Range("L17").Value = 9375.999899999999
$ws.Range("N17").Value = -9711.999899999999
$ws.Range("M17").Value = -5532
$ws.Range("H18").Value = 283.33334
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H58").Value = 83.666664
$ws.Range("I58").Value = 125
$ws.Range("J58").Value = 1
$ws.Range("K58").Value = 375
$ws.Range("L58").Value = 3
$ws.Range("M58").Value = -225
$ws.Range("N58").Value = -303
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 9000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -7502
$ws.Range("N99").ClearContents()
$ws.Range("H115").Value = 1180
$ws.Range("I115").Value = 1180
$ws.Range("K115").Value = 3540
$ws.Range("M115").Value = -1973
$ws.Range("H129").Value = 6000
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H136").Value = 78000
$ws.Range("J136").Value = 78000
$ws.Range("L136").Value = 78000
$ws.Range("N136").Value = -88200
$ws.Range("H138").Value = 1897.5
$ws.Range("J138").Value = 2500
$ws.Range("L138").Value = 7500
$ws.Range("N138").Value = -17780
$ws.Range("H140").Value = 122926.664
$ws.Range("J140").Value = 122926.664
$ws.Range("L140").Value = 122926.664
$ws.Range("N140").Value = -133286.664

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 4982
$ws.Range("J17").Value = 4982
$ws.Range("L17").Value = 4982
$ws.Range("N17").Value = -5328
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H102").Value = 1200
$ws.Range("I102").Value = 1200
$ws.Range("K102").Value = 1200
$ws.Range("M102").Value = 422
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 20
$ws.Range("I37").Value = 20
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 20
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 117
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 9300
$ws.Range("I56").Value = 9300
$ws.Range("K56").Value = 9300
$ws.Range("M56").Value = -8455
$ws.Range("H86").Value = 11247
$ws.Range("I86").Value = 9995
$ws.Range("K86").Value = 9995
$ws.Range("M86").Value = -8872
$ws.Range("H89").Value = 11247
$ws.Range("I89").Value = 9995
$ws.Range("K89").Value = 49975
$ws.Range("M89").Value = -44359
$ws.Range("H140").Value = 115445
$ws.Range("J140").Value = 115445
$ws.Range("L140").Value = 115445
$ws.Range("N140").Value = -125805

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 50
$ws.Range("J17").Value = 50
$ws.Range("L17").Value = 150
$ws.Range("N17").Value = -488
$ws.Range("H29").Value = 10
$ws.Range("J29").Value = 10
$ws.Range("L29").Value = 30
$ws.Range("N29").Value = -584
$ws.Range("H57").Value = 100
$ws.Range("J57").Value = 100
$ws.Range("L57").Value = 300
$ws.Range("N57").Value = -1418
$ws.Range("H58").Value = 1733.3334
$ws.Range("J58").Value = 4000
$ws.Range("L58").Value = 12000
$ws.Range("N58").Value = -12256
$ws.Range("H59").Value = 136.33333
$ws.Range("I59").Value = 136.33333
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 408.99999
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 131.00001
$ws.Range("N59").ClearContents()
$ws.Range("H68").Value = 1750
$ws.Range("J68").Value = 2500
$ws.Range("L68").Value = 7500
$ws.Range("N68").Value = -9122
$ws.Range("H71").Value = 1750
$ws.Range("J71").Value = 2500
$ws.Range("L71").Value = 22500
$ws.Range("N71").Value = -30612
$ws.Range("H81").Value = 125
$ws.Range("J81").Value = 100
$ws.Range("L81").Value = 300
$ws.Range("N81").Value = -2546
$ws.Range("H84").Value = 125
$ws.Range("J84").Value = 100
$ws.Range("L84").Value = 900
$ws.Range("N84").Value = -12132
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("H109").Value = 1388
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H133").Value = 930
$ws.Range("I133").Value = 930
$ws.Range("K133").Value = 2790
$ws.Range("M133").Value = 2270
$ws.Range("H134").Value = 7069.5713
$ws.Range("I134").Value = 7124.75
$ws.Range("K134").Value = 21374.25
$ws.Range("M134").Value = -16304.25
$ws.Range("H139").Value = 1223
$ws.Range("I139").Value = 1223
$ws.Range("K139").Value = 3669
$ws.Range("M139").Value = 1471
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 54.1
$ws.Range("I2").Value = 60.142857
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 60.142857
$ws.Range("L2").Value = 40
$ws.Range("M2").Value = 52.857143
$ws.Range("N2").Value = -266
$ws.Range("H82").Value = 1
$ws.Range("J82").Value = 1
$ws.Range("L82").Value = 1
$ws.Range("N82").Value = -767
$ws.Range("H85").Value = 1
$ws.Range("J85").Value = 1
$ws.Range("L85").Value = 1
$ws.Range("N85").Value = -2653
$ws.Range("H122").Value = 4002.3333
$ws.Range("I122").Value = 3003.5
$ws.Range("K122").Value = 9010.5
$ws.Range("M122").Value = -6560.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 448.66666
$ws.Range("I9").Value = 349
$ws.Range("J9").Value = 947
$ws.Range("K9").Value = 349
$ws.Range("L9").Value = 947
$ws.Range("M9").Value = -125
$ws.Range("N9").Value = -1395
$ws.Range("H19").Value = 2900
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2900
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2900
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -3240
$ws.Range("H35").Value = 899.6667
$ws.Range("I35").Value = 899.6667
$ws.Range("K35").Value = 899.6667
$ws.Range("M35").Value = -563.6667
$ws.Range("H58").Value = 8301.5
$ws.Range("I58").Value = 6500
$ws.Range("K58").Value = 6500
$ws.Range("M58").Value = -6240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1287
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H23").Value = 1201.3334
$ws.Range("I23").Value = 1201.3334
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1201.3334
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -972.3334
$ws.Range("N23").ClearContents()
